$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 2.056371915657621
$arr[1,0] = 1.962011255563993
$arr[2,0] = 1.906062748334932
$arr[3,0] = 1.88375973595538
$arr[4,0] = 1.88008621434011
$arr[5,0] = 1.905759956867087
$arr[6,0] = 2.023421140407663
$arr[7,0] = 2.270128714748239
$arr[8,0] = 2.461415179827668
$arr[9,0] = 2.550681091079014
$arr[10,0] = 2.584812294743244
$arr[11,0] = 2.577446857582288
$arr[12,0] = 2.553482481443666
$arr[13,0] = 2.538846486920647
$arr[14,0] = 2.455627097708998
$arr[15,0] = 2.40515400405252
$arr[16,0] = 2.376334582951415
$arr[17,0] = 2.366612981634091
$arr[18,0] = 2.410505041206079
$arr[19,0] = 2.560512455188643
$arr[20,0] = 2.660465948449655
$arr[21,0] = 2.606942050775558
$arr[22,0] = 2.408085220217288
$arr[23,0] = 2.201647925993029
$ws.Range("B2:B25").Value2 = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.122411137442441
$arr[1,0] = 0.1063295962782718
$arr[2,0] = 0.09646448688360465
$arr[3,0] = 0.09244584075437956
$arr[4,0] = 0.09177861299838241
$arr[5,0] = 0.096410285131725
$arr[6,0] = 0.1168638578462549
$arr[7,0] = 0.1570833747553877
$arr[8,0] = 0.1867565640766315
$arr[9,0] = 0.2002953163007533
$arr[10,0] = 0.2054288920741101
$arr[11,0] = 0.2043229693009323
$arr[12,0] = 0.2007175179606691
$arr[13,0] = 0.1985099818973026
$arr[14,0] = 0.1858726660707077
$arr[15,0] = 0.1781310671123038
$arr[16,0] = 0.1736820360677598
$arr[17,0] = 0.1721762828988744
$arr[18,0] = 0.1789547804663698
$arr[19,0] = 0.2017763356671196
$arr[20,0] = 0.2167313107301823
$arr[21,0] = 0.2087455989701823
$arr[22,0] = 0.1785823743815911
$arr[23,0] = 0.146185348816573
$ws.Range("C2:C25").Value2 = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.2119726927694643
$arr[1,0] = 0.2018537889925938
$arr[2,0] = 0.1956516979420542
$arr[3,0] = 0.1931264321116828
$arr[4,0] = 0.1927072246958943
$arr[5,0] = 0.1956176334867479
$arr[6,0] = 0.2084810446218199
$arr[7,0] = 0.2338221225137573
$arr[8,0] = 0.2525527555994529
$arr[9,0] = 0.261107458571189
$arr[10,0] = 0.264352525019433
$arr[11,0] = 0.2636533823620937
$arr[12,0] = 0.261374316197049
$arr[13,0] = 0.2599790704979057
$arr[14,0] = 0.2519944295173389
$arr[15,0] = 0.2471053188591554
$arr[16,0] = 0.2442963942540359
$arr[17,0] = 0.2433458629290328
$arr[18,0] = 0.247625440751392
$arr[19,0] = 0.2620435760207727
$arr[20,0] = 0.2714995482415645
$arr[21,0] = 0.2664494780669031
$arr[22,0] = 0.2473902878082868
$arr[23,0] = 0.2269499776279673
$ws.Range("D2:D25").Value2 = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 3.725823492963798
$arr[1,0] = 3.523603966045528
$arr[2,0] = 3.400570806908348
$arr[3,0] = 3.350708272844258
$arr[4,0] = 3.342444930876979
$arr[5,0] = 3.399897244842407
$arr[6,0] = 3.655857597994668
$arr[7,0] = 4.167233623051686
$arr[8,0] = 4.549385202219582
$arr[9,0] = 4.724788353229712
$arr[10,0] = 4.791445307841855
$arr[11,0] = 4.777078881858131
$arr[12,0] = 4.730267474656443
$arr[13,0] = 4.701625152144743
$arr[14,0] = 4.537954547029216
$arr[15,0] = 4.43795551157416
$arr[16,0] = 4.380584910167556
$arr[17,0] = 4.361184982211682
$arr[18,0] = 4.44858536811816
$arr[19,0] = 4.744010635710595
$arr[20,0] = 4.938466648969836
$arr[21,0] = 4.834551930613543
$arr[22,0] = 4.443779236496795
$arr[23,0] = 4.027811733591022
$ws.Range("F2:F25").Value2 = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.002526259695164789
$arr[1,0] = 0.002533270378531661
$arr[2,0] = 0.002537792721836755
$arr[3,0] = 0.002539690591864174
$arr[4,0] = 0.00254000905882336
$arr[5,0] = 0.002537818094299745
$arr[6,0] = 0.002528631918060086
$arr[7,0] = 0.002512335459138849
$arr[8,0] = 0.002501395230683662
$arr[9,0] = 0.002496639416310738
$arr[10,0] = 0.002494870047180697
$arr[11,0] = 0.002495249712579989
$arr[12,0] = 0.002496493218159302
$arr[13,0] = 0.002497259004242126
$arr[14,0] = 0.002501710461866509
$arr[15,0] = 0.002504497722658
$arr[16,0] = 0.00250612169124382
$arr[17,0] = 0.00250667511974097
$arr[18,0] = 0.002504198861710342
$arr[19,0] = 0.002496127115993865
$arr[20,0] = 0.002491035584401222
$arr[21,0] = 0.002493736282976472
$arr[22,0] = 0.00250433390973891
$arr[23,0] = 0.002516561681584822
$ws.Range("G2:G25").Value2 = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 0.3056714185922118
$arr[1,0] = 0.2950711110161279
$arr[2,0] = 0.2887651510484233
$arr[3,0] = 0.286245966390581
$arr[4,0] = 0.2858306984538643
$arr[5,0] = 0.2887309723327292
$arr[6,0] = 0.3019741098876239
$arr[7,0] = 0.3295728576130728
$arr[8,0] = 0.3508745985885895
$arr[9,0] = 0.360795086510123
$arr[10,0] = 0.364585375499928
$arr[11,0] = 0.3637675656267021
$arr[12,0] = 0.3611062391910309
$arr[13,0] = 0.3594804930076805
$arr[14,0] = 0.3502309472938521
$arr[15,0] = 0.344615974207386
$arr[16,0] = 0.3414080282078089
$arr[17,0] = 0.3403255738274993
$arr[18,0] = 0.3452114545716825
$arr[19,0] = 0.3618870188299326
$arr[20,0] = 0.3729816472690857
$arr[21,0] = 0.3670421124568719
$arr[22,0] = 0.3449421750914325
$arr[23,0] = 0.3219290194837185
$ws.Range("L2:L25").Value2 = $arr
